$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 196.54546
$ws.Range("I5").Value = 196.2
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 196.2
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -81.19999999999999
$ws.Range("N5").Value = -430
$ws.Range("H9").Value = 218.05405
$ws.Range("I9").Value = 289.88235
$ws.Range("K9").Value = 289.88235
$ws.Range("M9").Value = -120.88235
$ws.Range("H64").Value = 4160.7144
$ws.Range("I64").Value = 3412.5
$ws.Range("K64").Value = 3412.5
$ws.Range("M64").Value = -3164.5
$ws.Range("H67").Value = 4160.7144
$ws.Range("I67").Value = 3412.5
$ws.Range("K67").Value = 3412.5
$ws.Range("M67").Value = -2554.5
$ws.Range("H135").Value = 466.42426
$ws.Range("I135").Value = 313
$ws.Range("K135").Value = 2817
$ws.Range("M135").Value = -282
$ws.Range("H138").Value = 1515.4348
$ws.Range("I138").Value = 872.59375
$ws.Range("J138").Value = 2984.7856
$ws.Range("K138").Value = 2617.78125
$ws.Range("L138").Value = 8954.356800000001
$ws.Range("M138").Value = 2522.21875
$ws.Range("N138").Value = -19234.3568
$ws.Range("H141").Value = 45155.477
$ws.Range("I141").Value = 54642.945
$ws.Range("J141").Value = 11000.6
$ws.Range("K141").Value = 163928.835
$ws.Range("L141").Value = 33001.8
$ws.Range("M141").Value = -158748.835
$ws.Range("N141").Value = -43361.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4832.1113
$ws.Range("I2").Value = 4686.75
$ws.Range("K2").Value = 4686.75
$ws.Range("M2").Value = -4573.75
$ws.Range("H4").Value = 500875.25
$ws.Range("I4").Value = 667500.3
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 667500.3
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -667384.3
$ws.Range("N4").Value = -1232
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 1500
$ws.Range("K6").Value = 1500
$ws.Range("M6").Value = -1327
$ws.Range("H32").Value = 20401.453
$ws.Range("I32").Value = 3924.9692
$ws.Range("K32").Value = 3924.9692
$ws.Range("M32").Value = -3637.9692
$ws.Range("H45").Value = 483254.94
$ws.Range("I45").Value = 723087.9
$ws.Range("K45").Value = 723087.9
$ws.Range("M45").Value = -722710.9
$ws.Range("H61").Value = 1409.2778
$ws.Range("I61").Value = 1326.8667
$ws.Range("J61").Value = 1821.3334
$ws.Range("K61").Value = 1326.8667
$ws.Range("L61").Value = 1821.3334
$ws.Range("M61").Value = -1114.8667
$ws.Range("N61").Value = -2245.3334
$ws.Range("H97").Value = 705.3333
$ws.Range("I97").Value = 656.3889
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 656.3889
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = -160.3889
$ws.Range("N97").Value = -1991
$ws.Range("H116").Value = 4832.1113
$ws.Range("I116").Value = 4686.75
$ws.Range("K116").Value = 4686.75
$ws.Range("M116").Value = -2392.75
$ws.Range("H132").Value = 1273.3611
$ws.Range("I132").Value = 1042.0625
$ws.Range("J132").Value = 3123.75
$ws.Range("K132").Value = 3126.1875
$ws.Range("L132").Value = 9371.25
$ws.Range("M132").Value = -596.1875
$ws.Range("N132").Value = -14431.25
$ws.Range("H136").Value = 1409.2778
$ws.Range("I136").Value = 1326.8667
$ws.Range("J136").Value = 1821.3334
$ws.Range("K136").Value = 3980.6001
$ws.Range("L136").Value = 5464.0002
$ws.Range("M136").Value = -1430.6001
$ws.Range("N136").Value = -10564.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4832.1113
$ws.Range("I3").Value = 4686.75
$ws.Range("K3").Value = 4686.75
$ws.Range("M3").Value = -4572.75
$ws.Range("H94").Value = 1194.3422
$ws.Range("I94").Value = 1084.8823
$ws.Range("K94").Value = 1084.8823
$ws.Range("M94").Value = -633.8823
$ws.Range("H107").Value = 1880.2894
$ws.Range("I107").Value = 1220.1904
$ws.Range("K107").Value = 1220.1904
$ws.Range("M107").Value = 699.8096
$ws.Range("H134").Value = 973.1
$ws.Range("I134").Value = 970.6667
$ws.Range("J134").Value = 995
$ws.Range("K134").Value = 2912.0001
$ws.Range("L134").Value = 2985
$ws.Range("M134").Value = -377.0001000000002
$ws.Range("N134").Value = -8055

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 470.91666
$ws.Range("I22").Value = 633.3333
$ws.Range("J22").Value = 308.5
$ws.Range("K22").Value = 633.3333
$ws.Range("L22").Value = 308.5
$ws.Range("M22").Value = -283.3333
$ws.Range("N22").Value = -1008.5
$ws.Range("H31").Value = 35226.367
$ws.Range("I31").Value = 38651.594
$ws.Range("K31").Value = 38651.594
$ws.Range("M31").Value = -38356.594
$ws.Range("H34").Value = 35226.367
$ws.Range("I34").Value = 38651.594
$ws.Range("K34").Value = 38651.594
$ws.Range("M34").Value = -38449.594
$ws.Range("H107").Value = 1262.871
$ws.Range("I107").Value = 1102.0526
$ws.Range("K107").Value = 1102.0526
$ws.Range("M107").Value = 817.9474
$ws.Range("H134").Value = 3609.25
$ws.Range("I134").Value = 3548.8262
$ws.Range("K134").Value = 10646.4786
$ws.Range("M134").Value = -8111.4786

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 115.125
$ws.Range("I12").Value = 10.75
$ws.Range("J12").Value = 219.5
$ws.Range("K12").Value = 32.25
$ws.Range("L12").Value = 658.5
$ws.Range("M12").Value = 140.75
$ws.Range("N12").Value = -1004.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 284.10345
$ws.Range("I2").Value = 286.5263
$ws.Range("J2").Value = 279.5
$ws.Range("K2").Value = 286.5263
$ws.Range("L2").Value = 279.5
$ws.Range("M2").Value = -173.5263
$ws.Range("N2").Value = -505.5
$ws.Range("H122").Value = 1784.15
$ws.Range("J122").Value = 2936.25
$ws.Range("L122").Value = 8808.75
$ws.Range("N122").Value = -13708.75
$ws.Range("H126").Value = 3034.6365
$ws.Range("I126").Value = 2599.6667
$ws.Range("K126").Value = 7799.000100000001
$ws.Range("M126").Value = -5329.000100000001
$ws.Range("H132").Value = 20813.188
$ws.Range("I132").Value = 25193.615
$ws.Range("J132").Value = 1831.3334
$ws.Range("K132").Value = 75580.845
$ws.Range("L132").Value = 5494.0002
$ws.Range("M132").Value = -73050.845
$ws.Range("N132").Value = -10554.0002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H68").Value = 2586.375
$ws.Range("I68").Value = 1914
$ws.Range("K68").Value = 1914
$ws.Range("M68").Value = -1165
$ws.Range("H71").Value = 2586.375
$ws.Range("I71").Value = 1914
$ws.Range("K71").Value = 9570
$ws.Range("M71").Value = -5826
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H132").Value = 2675.2
$ws.Range("I132").Value = 2277.1372
$ws.Range("K132").Value = 6831.4116
$ws.Range("M132").Value = -4301.4116
$ws.Range("H136").Value = 2774.5796
$ws.Range("I136").Value = 2356.6924
$ws.Range("K136").Value = 7070.0772
$ws.Range("M136").Value = -4520.0772

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 11000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 11000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 11000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -11470
$ws.Range("H31").Value = 14845
$ws.Range("J31").Value = 14845
$ws.Range("L31").Value = 14845
$ws.Range("N31").Value = -15541
$ws.Range("H35").Value = 11000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 11000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 11000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -11580
$ws.Range("H62").Value = 38387.05
$ws.Range("I62").Value = 10516.833
$ws.Range("J62").Value = 50331.43
$ws.Range("K62").Value = 10516.833
$ws.Range("L62").Value = 50331.43
$ws.Range("M62").Value = -9892.833000000001
$ws.Range("N62").Value = -51579.43
$ws.Range("H65").Value = 38387.05
$ws.Range("I65").Value = 10516.833
$ws.Range("J65").Value = 50331.43
$ws.Range("K65").Value = 52584.165
$ws.Range("L65").Value = 251657.15
$ws.Range("M65").Value = -49464.165
$ws.Range("N65").Value = -257897.15
$ws.Range("H70").Value = 20090.625
$ws.Range("J70").Value = 20090.625
$ws.Range("L70").Value = 20090.625
$ws.Range("N70").Value = -20720.625
$ws.Range("H73").Value = 20090.625
$ws.Range("J73").Value = 20090.625
$ws.Range("L73").Value = 20090.625
$ws.Range("N73").Value = -22274.625
